$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C4 text: "16-64 years total" -> "16-64 years"
$ws.Range("C4").Value() = "16-64 years"

# Copy formats from column R into new column S (header + data rows)
$ws.Range("R3:R119").Copy()
$ws.Range("S3:S119").PasteSpecial(-4122)

# New year header
$ws.Range("S3").Value() = "2021"

# 2021 employment-rate data
$s2021 = @{
    4 = 77.1
    5 = 77.7
    6 = 76.7
    7 = 76.1
    8 = 79.1
    9 = 82.3
    10 = 79.7
    11 = 77.4
    12 = 73.6
    13 = 77
    14 = 73.2
    15 = 78.2
    16 = 79.1
    17 = 78.8
    18 = 78.9
    19 = 77.7
    20 = 77
    21 = 73.7
    22 = 77.1
    23 = 79
    24 = 82.1
    25 = 80.4
    26 = 84.3
    27 = 84.2
    28 = 79.4
    29 = 80.1
    30 = 81.9
    31 = 79.3
    32 = 77.3
    33 = 78
    34 = 80.8
    35 = 81
    36 = 79.9
    37 = 73.8
    38 = 73.7
    39 = 76
    40 = 76.9
    41 = 80.4
    42 = 81.1
    43 = 78.8
    44 = 81.6
    45 = 79.8
    46 = 83.9
    47 = 75.4
    48 = 78.3
    49 = 73.2
    50 = 78
    51 = 73.8
    52 = 67
    53 = 76.8
    54 = 73.4
    55 = 79
    56 = 74.7
    57 = 78
    58 = 79.2
    59 = 73.6
    60 = 76.5
    61 = 74.2
    62 = 76.7
    63 = 76.3
    64 = 77.2
    65 = 70.4
    66 = 80.2
    67 = 76.4
    68 = 73.9
    69 = 72.2
    70 = 73.5
    71 = 71.6
    72 = 78
    73 = 79.7
    74 = 77.7
    75 = 77.8
    76 = 76.3
    77 = 76.6
    78 = 79.3
    79 = 75.5
    80 = 75.8
    81 = 80.4
    82 = 81.1
    83 = 79.5
    84 = 75.5
    85 = 77.6
    86 = 77
    87 = 82.6
    88 = 82.2
    89 = 78.2
    90 = 74.4
    91 = 80
    92 = 76.4
    93 = 76.2
    94 = 80.3
    95 = 82.8
    96 = 78.5
    97 = 73.7
    98 = 79.1
    99 = 79.7
    100 = 80.1
    101 = 78.2
    102 = 79.4
    103 = 80.4
    104 = 77.9
    105 = 76.8
    106 = 78.3
    107 = 76.1
    108 = 76.1
    109 = 78.2
    110 = 75.2
    111 = 76.3
    112 = 77.6
    113 = 74.3
    114 = 77.3
    115 = 76.2
    116 = 82.5
    117 = 77.9
    118 = 77.1
    119 = 74.3
}
foreach ($row in $s2021.Keys) {
    $ws.Cells.Item([int]$row, 19).Value() = $s2021[$row]
}
